$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 3.281109666666667
$ws.Range("N2").Value = 9.843329000000001
$ws.Range("O2").Value = 0.2779739143628921
$ws.Range("P2").Value = 0.2779739143628921
$ws.Range("Q2").Value = 1.141817414374222
$ws.Range("R2").Value = 10.276356729368
$ws.Range("S2").Value = 0.2779739143628921
$ws.Range("T2").Value = 0.2779739143628921

# Row 3 updates
$ws.Range("M3").Value = 6.153936333333334
$ws.Range("O3").Value = 0.5213583040808726
$ws.Range("P3").Value = 0.5213583040808725
$ws.Range("S3").Value = 0.5213583040808726
$ws.Range("T3").Value = 0.5213583040808725

# Row 4 updates
$ws.Range("O4").Value = 0.2006677815562353
$ws.Range("P4").Value = 0.2006677815562353
$ws.Range("S4").Value = 0.2006677815562353
$ws.Range("T4").Value = 0.2006677815562353
